# --- vim_lib.xlsx edit: add vim entries on help, vimrc, fold, plugin etc ---

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")   # physically xl/worksheets/sheet1.xml, the tabSelected sheet

# Column C ("Description") now wraps text everywhere - apply this first (while the
# sheet still ends at row 14) so it picks up the header + every existing description
# cell, without materialising blank cells in the as-yet-unwritten new rows below.
$ws.Columns.Item(3).WrapText = $true

# Row 14 ("Insert <Tab> when expandtab ON") grows a little to fit the now-wrapped text.
$ws.Rows.Item(14).RowHeight = 30

# Row 15
$ws.Cells.Item(15, 1).Value = 'vim'
$ws.Cells.Item(15, 2).Value = 'help'

# Row 16
$ws.Cells.Item(16, 1).Value = 'vim'
$ws.Cells.Item(16, 2).Value = 'options'
$ws.Cells.Item(16, 3).Value = ':help ''option''
In the help text, point the cursor to tag ''option'' and go to tag with `Ctrl ]` to view all available options'
$ws.Cells.Item(16, 3).WrapText = $true
$ws.Rows.Item(16).RowHeight = 45

# Row 17
$ws.Cells.Item(17, 1).Value = 'vim'
$ws.Cells.Item(17, 2).Value = 'search'
$ws.Cells.Item(17, 3).Value = 'The incsearch option allows sync search before the enter is press when type the /{keyowrd} command:
:set incsearch'
$ws.Cells.Item(17, 3).WrapText = $true
$ws.Rows.Item(17).RowHeight = 45

# Row 18
$ws.Cells.Item(18, 1).Value = 'vim'
$ws.Cells.Item(18, 2).Value = 'search'
$ws.Cells.Item(18, 3).Value = '/{keyword} supports regular expression'
$ws.Cells.Item(18, 3).WrapText = $true

# Row 19
$ws.Cells.Item(19, 1).Value = 'vim'
$ws.Cells.Item(19, 2).Value = 'vimrc override hierarchy'
$ws.Cells.Item(19, 3).Value = '~/.vimrc (personal config) overrides /etc/vimrc (general config)'
$ws.Cells.Item(19, 3).WrapText = $true

# Row 20
$ws.Cells.Item(20, 1).Value = 'vim'
$ws.Cells.Item(20, 2).Value = 'fold'
$ws.Cells.Item(20, 3).Value = 'Set a fold:
zf    (zf and then use search or visual before zf to determine the endpoint of folding)
Open a fold:
zo   (or l at the beginning of line)
Close a fold:
zc'
$ws.Cells.Item(20, 3).WrapText = $true
$ws.Rows.Item(20).RowHeight = 90

# Row 21
$ws.Cells.Item(21, 1).Value = 'vim'
$ws.Cells.Item(21, 2).Value = 'error fix'
$ws.Cells.Item(21, 3).Value = 'When use cygwin to run vim plugins, there are many error with ^M, it is because of the unix file format (end of line) is not align with windows''. The solution is change the target script''s file format with vim with :set fileformat=unix'
$ws.Cells.Item(21, 3).WrapText = $true
$ws.Rows.Item(21).RowHeight = 45

# Row 22
$ws.Cells.Item(22, 1).Value = 'vim '
$ws.Cells.Item(22, 2).Value = 'edit multiple files'
$ws.Cells.Item(22, 3).Value = 'There are 3 ways to edit multiple files:
1. Use split open
2. Use tab open
3. vim *.java
For the vim *.java case, use :next / :n to move to next file'
$ws.Cells.Item(22, 3).WrapText = $true
$ws.Rows.Item(22).RowHeight = 75

# Row 23
$ws.Cells.Item(23, 1).Value = 'vim'
$ws.Cells.Item(23, 2).Value = 'plugin manager- pathogen'
$ws.Cells.Item(23, 3).Value = 'The plugin manager pathogen:
1. Installation: download the pathogen.vim to ~/.vim/autoload/pathogen.vim
2. Usage: download any plugin into ~/.vim/bundle/ and pathogen will auto-deplay the plugin, so no need to change config file'
$ws.Cells.Item(23, 3).WrapText = $true
$ws.Rows.Item(23).RowHeight = 60


# --- View state on the second worksheet ("Sheet1" tab / xl/worksheets/sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 31
$ws2.Range("A35:C42").Select()

# --- Restore the originally active sheet/selection ("Sheet2" tab / xl/worksheets/sheet1.xml) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("B24").Select()
